# LOQ4251.xlsx edit
# The original sheet had an orphan row 13 (B13/C13 = "198273 - Domingos Savio
# Giordani") sitting above "Programa resumido:" with no label in column A.
# That row is removed (shifting rows 14-22 up to 13-21), and several of the
# long-form paragraph cells are replaced with short placeholder values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 ("Objetivos:") loses its long descriptive paragraph, replaced with
# the docente id/name string (reused again further down, in "Método:").
$ws.Range("B10").Value = "198273 - Domingos Savio Giordani"
$ws.Range("C10").Value = "198273 - Domingos Savio Giordani"

# Remove the stray unlabeled row that held the "198273 - Domingos Savio
# Giordani" text above "Programa resumido:". This shifts rows 14:22 up to
# become rows 13:21 (labels/heights already line up after the shift).
$ws.Rows.Item(13).Delete()

# Row 13 ("Programa resumido:") content shortened.
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 ("Programa:") content replaced.
$ws.Range("B15").Value = "01/01/2021"
$ws.Range("C15").Value = "01/01/2021"

# Row 18 ("Método:") now carries the docente info that used to sit in the
# removed orphan row.
$ws.Range("B18").Value = "198273 - Domingos Savio Giordani"
$ws.Range("C18").Value = "198273 - Domingos Savio Giordani"

# Row 19 ("Critério:") picks up the text that used to belong to "Método:".
$ws.Range("B19").Value = "Duas provas escritas e um seminário que, juntos, constituem a primeira avaliação."
$ws.Range("C19").Value = "Duas provas escritas e um seminário que, juntos, constituem a primeira avaliação."

# Row 20 ("Norma de recuperação:") picks up the text that used to belong to
# "Critério:".
$ws.Range("B20").Value = "A nota de primeira avaliação será igual à média das notas das duas provas, com peso 7 somada à nota do seminário com peso 3. Alunos com nota de primeira avaliação igual ou superior a 5 estarão aprovados, com nota entre 3 e 4,9 em recuperação e abaixo de 3 reprovados."
$ws.Range("C20").Value = "A nota de primeira avaliação será igual à média das notas das duas provas, com peso 7 somada à nota do seminário com peso 3. Alunos com nota de primeira avaliação igual ou superior a 5 estarão aprovados, com nota entre 3 e 4,9 em recuperação e abaixo de 3 reprovados."

# Row 21 ("Bibliografia:") picks up the text that used to belong to "Norma
# de recuperação:"; the old bibliography paragraph is dropped entirely.
$ws.Range("B21").Value = "A recuperação se constituirá de uma prova abordando todos os assuntos do semestre, a nota de segunda avaliação será igual à média entre a nota de primeira avaliação e a prova de recuperação. Alunos com nota de segunda avaliação igual ou superior a 5 estarão aprovados e inferior a 5 reprovados."
$ws.Range("C21").Value = "A recuperação se constituirá de uma prova abordando todos os assuntos do semestre, a nota de segunda avaliação será igual à média entre a nota de primeira avaliação e a prova de recuperação. Alunos com nota de segunda avaliação igual ou superior a 5 estarão aprovados e inferior a 5 reprovados."
